$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# --- Row 12 fixes -----------------------------------------------------
# Apply the "normal data" style (same as sibling rows) to the cells that were
# missing it, and correct the run-time value in B12.
$ws.Range("A12").Style = $ws.Range("A11").Style
$ws.Range("C12").Style = $ws.Range("C11").Style
$ws.Range("E12").Style = $ws.Range("E11").Style
$ws.Range("F12").Style = $ws.Range("F11").Style
$ws.Range("G12").Style = $ws.Range("G11").Style

$ws.Range("B12").Value = 44363.80596668981

# --- Row 13: new interview-history entry, fully styled like row 11/12 -
$ws.Range("A13").Style = $ws.Range("A11").Style
$ws.Range("A13").Value = "2021-06-17"

$ws.Range("B13").Style = $ws.Range("B11").Style
$ws.Range("B13").Value = 44364.57609354167

$ws.Range("C13").Style = $ws.Range("C11").Style
$ws.Range("C13").Value = "145_data_hstry"

$ws.Range("D13").Style = $ws.Range("D11").Style
$ws.Range("D13").Value = 165

$ws.Range("E13").Style = $ws.Range("E11").Style
$ws.Range("E13").Value = 164

$ws.Range("F13").Style = $ws.Range("F11").Style
$ws.Range("F13").Value = 1

$ws.Range("G13").Style = $ws.Range("G11").Style
$ws.Range("G13").Value = 4.91

# --- Row 14: new interview-history entry, default style (like row 1-9) -
$ws.Range("A14").Value = "2021-06-17"

$ws.Range("B14").Style = $ws.Range("B11").Style
$ws.Range("B14").Value = 44364.60561181853

$ws.Range("C14").Value = "145_hstry_data"

$ws.Range("D14").Style = $ws.Range("D11").Style
$ws.Range("D14").Value = 165

$ws.Range("E14").Value = 164
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.95
